# PM-Sheet_BAUX(1).xlsx edit
# ---------------------------------------------------------------
# Re-uploaded copy of the PM sheet ("Add files via upload"): the
# worksheet title becomes "PM-Sheet (1)", the deadline moves a week
# later, the "Responsible" column (I) gets filled in for every task
# row, the "Remain" value for the in-progress task is set to 0, and
# the active selection ends up on G13 with the view scrolled back to
# the top.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell: "PM-Sheet" -> "PM-Sheet (1)"
$ws.Range("B1").Value = "PM-Sheet (1)"

# Deadline: 2019-01-24 -> 2019-01-31 (serial 43489 -> 43496)
$ws.Range("D4").Value = 43496

# Responsible column (I) now filled in for each task row. Pick up the
# (unstyled) formatting of the neighbouring Task cell first so the new
# cells don't pick up an explicit style of their own, matching how the
# rest of that column looks.
$ws.Range("C12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I12").Value = "Djukic, Hamzic and Taha"
$ws.Range("I13").Value = "Djukic, Hamzic and Taha"
$ws.Range("I14").Value = "Djukic, Hamzic and Taha"

# Remain column (H) for the "Make a login..." row now shows 0
$ws.Range("H13").Value = 0

# Move the selection to G13 (and let the view scroll back to the top)
[void]$ws.Range("G13").Select()
